# Insert two new price-report rows for "Apio" (Vega Modelo de Temuco) just
# above the existing row 441, shifting the old rows 441:547 down to 443:549.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything from row 441 down by inserting two blank rows above it.
$ws.Rows("441:442").Insert()

# New row 441 data
$ws.Range("A441").Value = 10
$ws.Range("B441").Value = "Vega Modelo de Temuco"
$ws.Range("C441").Value = "La Araucanía"
$ws.Range("D441").Value = 45173
$ws.Range("E441").Value = 9
$ws.Range("F441").Value = 100112017
$ws.Range("G441").Value = "Apio"
$ws.Range("H441").Value = "Americana (o)"
$ws.Range("I441").Value = "Primera"
$ws.Range("J441").Value = 200
$ws.Range("K441").Value = 8000
$ws.Range("L441").Value = 8000
$ws.Range("M441").Value = 8000
$ws.Range("N441").Value = '$/caja 8 unidades'
$ws.Range("O441").Value = "Provincia del Elquí"
$ws.Range("P441").Value = 8000
$ws.Range("Q441").Value = 1
$ws.Range("R441").Value = "Hortaliza"

# New row 442 data
$ws.Range("A442").Value = 10
$ws.Range("B442").Value = "Vega Modelo de Temuco"
$ws.Range("C442").Value = "La Araucanía"
$ws.Range("D442").Value = 45173
$ws.Range("E442").Value = 9
$ws.Range("F442").Value = 100112017
$ws.Range("G442").Value = "Apio"
$ws.Range("H442").Value = "Americana (o)"
$ws.Range("I442").Value = "Primera"
$ws.Range("J442").Value = 140
$ws.Range("K442").Value = 8000
$ws.Range("L442").Value = 8000
$ws.Range("M442").Value = 8000
$ws.Range("N442").Value = '$/docena de matas'
$ws.Range("O442").Value = "Provincia del Elquí"
$ws.Range("P442").Value = 1333
$ws.Range("Q442").Value = 6
$ws.Range("R442").Value = "Hortaliza"
